# "ajout temps par personne" — fill in the per-person hour counts and
# drop the now-unused "Code IHM" task column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# The "Code IHM" column (H) is removed entirely; Excel shifts
# Tests/Rapport/Gestion de projet (and the Total column) one slot to
# the left, keeps formulas consistent and shrinks the E2:H2 merge to
# E2:G2 automatically.
$ws.Columns.Item(8).Delete()

# --- fill in the hours-per-person grid (columns E:J, rows 6:12) ---
# D6 Nicolas NATIVEL
$ws.Cells.Item(6, 5).Value = 6
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(6, 7).Value = 40
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = $null
$ws.Cells.Item(6, 10).Value = 4

# D7 Romain MATHONAT
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 5
$ws.Cells.Item(7, 7).Value = 45
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0

# D8 Mathieu GAILLARD
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 10
$ws.Cells.Item(8, 7).Value = 18
$ws.Cells.Item(8, 8).Value = 5
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(8, 10).Value = 1

# D9 Mohammed EL ARASS, D10 Guillaume KHENG, D11 Thomas FAVROT,
# D12 Killian OLLIVIER stay untouched (no hours logged yet).

# Move the active selection the way the author left it.
$ws.Range("D3").Select()

Write-Output "done"
